$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Save" in H1, copying the header style (bold/border/centered)
# from the neighboring "sum" header (G1) before writing the text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill H2:H6 with 0 (plain data style, matching the other numeric columns)
$ws.Range("H2:H6").Value = 0
